$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8931353688240051
$ws.Range("B1").Value = 1.495593786239624
$ws.Range("C1").Value = 8.791707038879395
$ws.Range("D1").Value = 2.073096513748169
$ws.Range("E1").Value = 1.203125238418579
